$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 74
$ws1.Range("H74").Value = 3502
$ws1.Range("I74").Value = 3502
$ws1.Range("K74").Value = 3502
$ws1.Range("M74").Value = -2566
# row 77
$ws1.Range("H77").Value = 3502
$ws1.Range("I77").Value = 3502
$ws1.Range("K77").Value = 17510
$ws1.Range("M77").Value = -12830
# row 80
$ws1.Range("H80").Value = 6212.5557
$ws1.Range("J80").Value = 5462.273
$ws1.Range("L80").Value = 16386.819
$ws1.Range("N80").Value = -18382.819
# row 83
$ws1.Range("H83").Value = 6212.5557
$ws1.Range("J83").Value = 5462.273
$ws1.Range("L83").Value = 49160.457
$ws1.Range("N83").Value = -59144.457
# row 98
$ws1.Range("H98").Value = 1236.1
$ws1.Range("I98").Value = 1111.6666
$ws1.Range("K98").Value = 1111.6666
$ws1.Range("M98").Value = 386.3334
# row 112
$ws1.Range("H112").Value = 8755.1875
$ws1.Range("J112").Value = 8755.1875
$ws1.Range("L112").Value = 26265.5625
$ws1.Range("N112").Value = -28481.5625
# row 122
$ws1.Range("H122").Value = 1236.1
$ws1.Range("I122").Value = 1111.6666
$ws1.Range("K122").Value = 3334.9998
$ws1.Range("M122").Value = -884.9998000000001
# row 137
$ws1.Range("H137").Value = 20602.766
$ws1.Range("I137").Value = 32325.1
$ws1.Range("K137").Value = 96975.29999999999
$ws1.Range("M137").Value = -94425.29999999999

# --- ARM ---
# row 32
$ws2.Range("H32").Value = 19986.508
$ws2.Range("I32").Value = 20308.016
$ws2.Range("J32").Value = 15999.8
$ws2.Range("K32").Value = 20308.016
$ws2.Range("L32").Value = 15999.8
$ws2.Range("M32").Value = -20021.016
$ws2.Range("N32").Value = -16573.8
# row 43
$ws2.Range("H43").Value = 92377
$ws2.Range("J43").Value = 92377
$ws2.Range("L43").Value = 92377
$ws2.Range("N43").Value = -93003
# row 63
$ws2.Range("H63").Value = 4066.6
$ws2.Range("I63").Value = 1333
$ws2.Range("J63").Value = 4750
$ws2.Range("K63").Value = 1333
$ws2.Range("L63").Value = 4750
$ws2.Range("M63").Value = -647
$ws2.Range("N63").Value = -6122
# row 66
$ws2.Range("H66").Value = 4066.6
$ws2.Range("I66").Value = 1333
$ws2.Range("J66").Value = 4750
$ws2.Range("K66").Value = 6665
$ws2.Range("L66").Value = 23750
$ws2.Range("M66").Value = -3233
$ws2.Range("N66").Value = -30614
# row 74
$ws2.Range("H74").Value = 463778.53
$ws2.Range("I74").Value = 857946.3
$ws2.Range("J74").Value = 3916.1667
$ws2.Range("K74").Value = 857946.3
$ws2.Range("L74").Value = 3916.1667
$ws2.Range("M74").Value = -857072.3
$ws2.Range("N74").Value = -5664.1667
# row 77
$ws2.Range("H77").Value = 463778.53
$ws2.Range("I77").Value = 857946.3
$ws2.Range("J77").Value = 3916.1667
$ws2.Range("K77").Value = 4289731.5
$ws2.Range("L77").Value = 19580.8335
$ws2.Range("M77").Value = -4285363.5
$ws2.Range("N77").Value = -28316.8335
# row 110
$ws2.Range("H110").Value = 2102
$ws2.Range("I110").Value = 2102
$ws2.Range("K110").Value = 2102
$ws2.Range("M110").Value = -57
# row 132
$ws2.Range("H132").Value = 3202.5813
$ws2.Range("I132").Value = 1242.9166
$ws2.Range("K132").Value = 3728.7498
$ws2.Range("M132").Value = -1198.7498

# --- BSM ---
# row 86
$ws3.Range("H86").Value = 6492.385
$ws3.Range("I86").Value = 2415.1667
$ws3.Range("J86").Value = 9987.143
$ws3.Range("K86").Value = 2415.1667
$ws3.Range("L86").Value = 9987.143
$ws3.Range("M86").Value = -1292.1667
$ws3.Range("N86").Value = -12233.143
# row 89
$ws3.Range("H89").Value = 6492.385
$ws3.Range("I89").Value = 2415.1667
$ws3.Range("J89").Value = 9987.143
$ws3.Range("K89").Value = 12075.8335
$ws3.Range("L89").Value = 49935.715
$ws3.Range("M89").Value = -6459.833500000001
$ws3.Range("N89").Value = -61167.715
# row 134
$ws3.Range("H134").Value = 2179.724
$ws3.Range("I134").Value = 2035.0834
$ws3.Range("J134").Value = 2874
$ws3.Range("K134").Value = 6105.2502
$ws3.Range("L134").Value = 8622
$ws3.Range("M134").Value = -3570.2502
$ws3.Range("N134").Value = -13692

# --- CRP ---
# row 86
$ws4.Range("H86").Value = 10929.762
$ws4.Range("I86").Value = 5894.077
$ws4.Range("J86").Value = 19112.75
$ws4.Range("K86").Value = 5894.077
$ws4.Range("L86").Value = 19112.75
$ws4.Range("M86").Value = -4771.077
$ws4.Range("N86").Value = -21358.75
# row 89
$ws4.Range("H89").Value = 10929.762
$ws4.Range("I89").Value = 5894.077
$ws4.Range("J89").Value = 19112.75
$ws4.Range("K89").Value = 29470.385
$ws4.Range("L89").Value = 95563.75
$ws4.Range("M89").Value = -23854.385
$ws4.Range("N89").Value = -106795.75
# row 105
$ws4.Range("H105").Value = 2385.077
$ws4.Range("I105").Value = 1228.375
$ws4.Range("J105").Value = 4235.8
$ws4.Range("K105").Value = 1228.375
$ws4.Range("L105").Value = 4235.8
$ws4.Range("M105").Value = 518.625
$ws4.Range("N105").Value = -7729.8
# row 107
$ws4.Range("H107").Value = 423.1111
$ws4.Range("I107").Value = 330.35294
$ws4.Range("K107").Value = 330.35294
$ws4.Range("M107").Value = 1589.64706
# row 134
$ws4.Range("H134").Value = 3215.9333
$ws4.Range("I134").Value = 2855.818
$ws4.Range("K134").Value = 8567.454000000002
$ws4.Range("M134").Value = -6032.454000000002

# --- CUL ---
# row 57
$ws5.Range("H57").Value = 9699.76
$ws5.Range("J57").Value = 9699.76
$ws5.Range("L57").Value = 29099.28
$ws5.Range("N57").Value = -30217.28
# row 116
$ws5.Range("H116").Value = 8392.478999999999
$ws5.Range("I116").Value = 2765.4
$ws5.Range("J116").Value = 9955.556
$ws5.Range("K116").Value = 8296.200000000001
$ws5.Range("L116").Value = 29866.668
$ws5.Range("M116").Value = -4854.200000000001
$ws5.Range("N116").Value = -36750.66800000001
# row 131
$ws5.Range("H131").Value = 328701.53
$ws5.Range("I131").Value = 608574.3
$ws5.Range("J131").Value = 2183.3333
$ws5.Range("K131").Value = 1825722.9
$ws5.Range("L131").Value = 6549.999899999999
$ws5.Range("M131").Value = -1820682.9
$ws5.Range("N131").Value = -16629.9999
# row 136
$ws5.Range("H136").Value = 2384.6667
$ws5.Range("I136").Value = 2061.6
$ws5.Range("K136").Value = 6184.799999999999
$ws5.Range("M136").Value = -1084.799999999999
# row 140
$ws5.Range("H140").Value = 3292.5
$ws5.Range("I140").Value = 3292.5
$ws5.Range("K140").Value = 9877.5
$ws5.Range("M140").Value = -4697.5

# --- GSM ---
# row 107
$ws6.Range("H107").Value = 254.16667
$ws6.Range("I107").Value = 199.5
$ws6.Range("J107").Value = 527.5
$ws6.Range("K107").Value = 199.5
$ws6.Range("L107").Value = 527.5
$ws6.Range("M107").Value = 1720.5
$ws6.Range("N107").Value = -4367.5
# row 122
$ws6.Range("H122").Value = 3626.375
$ws6.Range("I122").Value = 3317.3
$ws6.Range("J122").Value = 4141.5
$ws6.Range("K122").Value = 9951.900000000001
$ws6.Range("L122").Value = 12424.5
$ws6.Range("M122").Value = -7501.900000000001
$ws6.Range("N122").Value = -17324.5
# row 126
$ws6.Range("H126").Value = 3679
$ws6.Range("I126").Value = 2798.3333
$ws6.Range("J126").Value = 5000
$ws6.Range("K126").Value = 8394.999899999999
$ws6.Range("L126").Value = 15000
$ws6.Range("M126").Value = -5924.999899999999
$ws6.Range("N126").Value = -19940
# row 136
$ws6.Range("H136").Value = 9872.666999999999
$ws6.Range("J136").Value = 9872.666999999999
$ws6.Range("L136").Value = 29618.001
$ws6.Range("N136").Value = -34718.001

# --- LTW ---
# row 14
$ws7.Range("H14").Value = 0
$ws7.Range("J14").Value = 0
$ws7.Range("L14").Value = 0
$ws7.Range("N14").ClearContents()
# row 40
$ws7.Range("H40").Value = 3793.8
$ws7.Range("I40").Value = 3548.7222
$ws7.Range("K40").Value = 3548.7222
$ws7.Range("M40").Value = -3412.7222
# row 46
$ws7.Range("H46").Value = 5139.1113
$ws7.Range("I46").Value = 1086.1428
$ws7.Range("J46").Value = 7718.273
$ws7.Range("K46").Value = 1086.1428
$ws7.Range("L46").Value = 7718.273
$ws7.Range("M46").Value = -898.1428000000001
$ws7.Range("N46").Value = -8094.273
# row 55
$ws7.Range("H55").Value = 1534.4445
$ws7.Range("I55").Value = 540
$ws7.Range("K55").Value = 540
$ws7.Range("M55").Value = -367
# row 61
$ws7.Range("H61").Value = 917.0833
$ws7.Range("I61").Value = 934.36365
$ws7.Range("J61").Value = 727
$ws7.Range("K61").Value = 934.36365
$ws7.Range("L61").Value = 727
$ws7.Range("M61").Value = -732.36365
$ws7.Range("N61").Value = -1131
# row 93
$ws7.Range("H93").Value = 2856.8518
$ws7.Range("I93").Value = 2140.0557
$ws7.Range("K93").Value = 2140.0557
$ws7.Range("M93").Value = -892.0556999999999
# row 99
$ws7.Range("H99").Value = 27142.25
$ws7.Range("I99").Value = 9999.5
$ws7.Range("K99").Value = 9999.5
$ws7.Range("M99").Value = -7004.5
# row 113
$ws7.Range("H113").Value = 917.0833
$ws7.Range("I113").Value = 934.36365
$ws7.Range("J113").Value = 727
$ws7.Range("K113").Value = 934.36365
$ws7.Range("L113").Value = 727
$ws7.Range("M113").Value = 1235.63635
$ws7.Range("N113").Value = -5067
# row 122
$ws7.Range("H122").Value = 8223.806
$ws7.Range("I122").Value = 8526.174000000001
$ws7.Range("K122").Value = 25578.522
$ws7.Range("M122").Value = -23128.522
# row 132
$ws7.Range("H132").Value = 4448.8125
$ws7.Range("I132").Value = 4545.467
$ws7.Range("K132").Value = 13636.401
$ws7.Range("M132").Value = -11106.401
# row 136
$ws7.Range("H136").Value = 5992.3125
$ws7.Range("I136").Value = 6262.3335
$ws7.Range("J136").Value = 5476.8184
$ws7.Range("K136").Value = 18787.0005
$ws7.Range("L136").Value = 16430.4552
$ws7.Range("M136").Value = -16237.0005
$ws7.Range("N136").Value = -21530.4552

# --- WVR ---
# row 58
$ws8.Range("H58").Value = 49950
$ws8.Range("I58").Value = 0
$ws8.Range("J58").Value = 49950
$ws8.Range("K58").Value = 0
$ws8.Range("L58").Value = 49950
$ws8.Range("N58").Value = -50566
$ws8.Range("M58").ClearContents()
# row 132
$ws8.Range("H132").Value = 4582.8335
$ws8.Range("J132").Value = 2722.25
$ws8.Range("L132").Value = 8166.75
$ws8.Range("N132").Value = -13226.75
# row 135
$ws8.Range("H135").Value = 99992.39999999999
$ws8.Range("J135").Value = 99992.39999999999
$ws8.Range("L135").Value = 99992.39999999999
$ws8.Range("N135").Value = -110132.4
# row 136
$ws8.Range("H136").Value = 26797.037
$ws8.Range("I136").Value = 29932.25
$ws8.Range("J136").Value = 1715.3334
$ws8.Range("K136").Value = 89796.75
$ws8.Range("L136").Value = 5146.0002
$ws8.Range("M136").Value = -87246.75
$ws8.Range("N136").Value = -10246.0002
